# ---------------------------------------------------------------------------
# Applies the "Added getSchedule tests ..." commit:
#  - Adjust selection on the Organizations sheet
#  - Deselect the Practitioner tab (it is no longer the active sheet)
#  - Add a brand new "Sheet4" worksheet after Practitioner describing the
#    fixed + rolling slot pre-requisites for the getSchedule tests, and make
#    it the active sheet/tab.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$patients      = $wb.Worksheets.Item(1)
$organizations = $wb.Worksheets.Item(2)
$practitioner  = $wb.Worksheets.Item(3)

# --- Organizations: move the saved selection from N15 to A3 ---------------
$organizations.Range("A3").Select()

# --- Practitioner: keep its own selection (D6) but it stops being the
#     active/selected tab once we activate the new sheet below -------------
$practitioner.Activate()
$practitioner.Range("D6").Select()

# --- Add the new worksheet right after "Practitioner" ----------------------
$sheet4 = $wb.Worksheets.Add($null, $practitioner)
$sheet4.Name = "Sheet4"

# Column widths (approximate Excel's internal pixel-quantised widths)
$sheet4.Columns.Item(1).ColumnWidth = 3.86
$sheet4.Columns.Item(2).ColumnWidth = 18.55
$sheet4.Columns.Item(3).ColumnWidth = 13.72
$sheet4.Columns.Item(4).ColumnWidth = 8.48

# Colours reused throughout the sheet
$headerFill = 14540253   # FFDDDDDD / FFD6DCE5 (existing light-grey fill)
$headerFillPattern = 15064278
$columnHeadFill = 13421772   # FFCCCCCC / FFBFBFBF (new medium-grey fill)
$columnHeadFillPattern = 12566463

# ---------------------------------------------------------------------------
# Section 1: "Required Fixed Slots for testing"
# ---------------------------------------------------------------------------
$sheet4.Range("B2:D2").Merge()
$sheet4.Range("B2").Value = "Required Fixed Slots for testing"
$sheet4.Range("B2:D2").Font.Size = 13
$sheet4.Range("B2:D2").Font.Bold = $true
$sheet4.Range("B2:D2").Interior.Color = $headerFill
$sheet4.Range("B2:D2").Interior.PatternColor = $headerFillPattern
$sheet4.Rows.Item(2).RowHeight = 16.15

$sheet4.Range("B3:D3").Interior.Color = $headerFill
$sheet4.Range("B3:D3").Interior.PatternColor = $headerFillPattern

$sheet4.Range("B4").Value = "Organization Code"
$sheet4.Range("C4").Value = "Date"
$sheet4.Range("D4").Value = "Time"
$sheet4.Range("B4:D4").Interior.Color = $columnHeadFill
$sheet4.Range("B4:D4").Interior.PatternColor = $columnHeadFillPattern

$sheet4.Range("B5").Value = "ORG1"
$sheet4.Range("C5").Value = 42792
$sheet4.Range("C5").NumberFormat = "YYYY\-MM\-DD"
$sheet4.Range("D5").Value = 0.440972222222222
$sheet4.Range("D5").NumberFormat = "HH:MM:SS"
$sheet4.Range("B5:D5").Interior.Color = $headerFill
$sheet4.Range("B5:D5").Interior.PatternColor = $headerFillPattern

$sheet4.Range("B6").Value = "ORG1"
$sheet4.Range("C6").Value = 42737
$sheet4.Range("C6").NumberFormat = "YYYY\-MM\-DD"
$sheet4.Range("D6").Value = 0.527777777777778
$sheet4.Range("D6").NumberFormat = "HH:MM:SS"
$sheet4.Range("B6:D6").Interior.Color = $headerFill
$sheet4.Range("B6:D6").Interior.PatternColor = $headerFillPattern

$sheet4.Range("B7").Value = "ORG1"
$sheet4.Range("C7").Value = 43098
$sheet4.Range("C7").NumberFormat = "YYYY\-MM\-DD"
$sheet4.Range("D7").Value = 0.527777777777778
$sheet4.Range("D7").NumberFormat = "HH:MM:SS"
$sheet4.Range("B7:D7").Interior.Color = $headerFill
$sheet4.Range("B7:D7").Interior.PatternColor = $headerFillPattern

# ---------------------------------------------------------------------------
# Section 2: "Rolling Slots"
# ---------------------------------------------------------------------------
$sheet4.Range("B10").Value = "Rolling Slots"
$sheet4.Range("B10:M10").Font.Size = 13
$sheet4.Range("B10:M10").Font.Bold = $true
$sheet4.Range("B10:M10").Interior.Color = $headerFill
$sheet4.Range("B10:M10").Interior.PatternColor = $headerFillPattern
$sheet4.Rows.Item(10).RowHeight = 16.15

$sheet4.Range("B11:M11").Interior.Color = $headerFill
$sheet4.Range("B11:M11").Interior.PatternColor = $headerFillPattern

$sheet4.Range("B12:M12").Merge()
$sheet4.Range("B12").Value = "A set of rolling/available slots needs to be present from the current date for up to a 14 days into the future with the following number of slots for each of the following organizations."
$sheet4.Range("B12:M12").Interior.Color = $headerFill
$sheet4.Range("B12:M12").Interior.PatternColor = $headerFillPattern

$sheet4.Range("B13:M13").Interior.Color = $headerFill
$sheet4.Range("B13:M13").Interior.PatternColor = $headerFillPattern

$sheet4.Range("B14").Value = "Organization Code"
$sheet4.Range("C14").Value = "Number of slots"
$sheet4.Range("B14:C14").Interior.Color = $columnHeadFill
$sheet4.Range("B14:C14").Interior.PatternColor = $columnHeadFillPattern
$sheet4.Range("D14:M14").Interior.Color = $headerFill
$sheet4.Range("D14:M14").Interior.PatternColor = $headerFillPattern

$sheet4.Range("B15").Value = "ORG1"
$sheet4.Range("C15").Value = 30
$sheet4.Range("B15:M15").Interior.Color = $headerFill
$sheet4.Range("B15:M15").Interior.PatternColor = $headerFillPattern

$sheet4.Range("B16").Value = "ORG2"
$sheet4.Range("C16").Value = 20
$sheet4.Range("B16:M16").Interior.Color = $headerFill
$sheet4.Range("B16:M16").Interior.PatternColor = $headerFillPattern

$sheet4.Range("B17").Value = "ORG3"
$sheet4.Range("C17").Value = 0
$sheet4.Range("B17:M17").Interior.Color = $headerFill
$sheet4.Range("B17:M17").Interior.PatternColor = $headerFillPattern

# --- Make Sheet4 the active sheet/tab, mirroring the saved selection ------
$sheet4.Activate()
$excel.ActiveWindow.Zoom = 85
$sheet4.Range("B19").Select()
